# Edit script: apply the commit "small edits to ppt"
#  1. Update the "Goal" slide (slide 2) body text (impact -> significance, etc.)
#  2. Insert a new "The Data" (Title and Content) slide at position 3.

$p = $ppt.ActivePresentation

# --- 1. Update body text on slide 2 ("Goal") -----------------------------
$goalSlide = $p.Slides.Item(2)
$goalBody = $goalSlide.Shapes.Item("Content Placeholder 2")
$goalBody.TextFrame.TextRange.Text = "Assess the significance various factors such as age, race, mental illness and gender have on the population of people in the United States who have been fatally shot by the police. "

# --- 2. Insert new "The Data" slide at position 3 ------------------------
# ppLayoutText (2) => "Title and Content" custom layout (Title + body placeholder)
$newSlide = $p.Slides.Add(3, 2)

$newSlide.Shapes.Item("Title 1").TextFrame.TextRange.Text = "The Data"

$body = $newSlide.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
$body.Text = "Fatal police shootings in the US from 2015-2021"
[void]$body.InsertAfter("`rCollected by the Washington Post")
[void]$body.InsertAfter("`rDetails about each individual: race, gender, age, location, signs of mental illness, ")
[void]$body.InsertAfter("etc")
